$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("D2").Value = 15.54
$ws.Range("E2").Value = 56
$ws.Range("F2").Value = 10.02
$ws.Range("K2").Value = 57.5
$ws.Range("N2").Value = 53.62998959737769

# Row 3 updates
$ws.Range("D3").Value = 8.66
$ws.Range("E3").Value = 58.7
$ws.Range("F3").Value = 15.69
$ws.Range("H3").Value = 36
$ws.Range("I3").Value = 50
$ws.Range("K3").Value = 54.1
$ws.Range("N3").Value = 53.62998959737769
